$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 87, shifting the existing
# rows 87-92 down to 88-93 (this also grows the used range to T93).
$ws.Rows(87).Insert()

# Populate the newly inserted row 87 with the new weekly price record.
# Columns A,B,C,E,F,G,H,I,J,K,L,Q,R,T follow the same pattern as the
# surrounding "Coco" rows; only D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado) and S (Precio $/Kg) change.
$ws.Range("A87").Value = 10
$ws.Range("B87").Value = "Vega Modelo de Temuco"
$ws.Range("C87").Value = "La Araucanía"
$ws.Range("D87").Value = 44783
$ws.Range("D87").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E87").Value = 9
$ws.Range("F87").Value = "Fruta"
$ws.Range("G87").Value = 100108
$ws.Range("H87").Value = "Tropicales y subtropicales"
$ws.Range("I87").Value = 100108007
$ws.Range("J87").Value = "Coco"
$ws.Range("K87").Value = "Sin especificar"
$ws.Range("L87").Value = "Primera"
$ws.Range("M87").Value = 35
$ws.Range("N87").Value = 30000
$ws.Range("O87").Value = 30000
$ws.Range("P87").Value = 30000
$ws.Range("Q87").Value = "$/malla 20 unidades"
$ws.Range("R87").Value = "Perú"
$ws.Range("S87").Value = 1500
$ws.Range("T87").Value = 20
